$d = $word.ActiveDocument

# 1) "Unity Version: 2019.3" -> "Unity Version: 2019.3.0f5"
#    The new text is appended as its own run right after the existing
#    "2019.3" run (matches a user simply typing more text at that spot).
$r1 = $d.Content
$r1.Find.Execute("2019.3", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $d.Range($r1.End, $r1.End)
$insertPoint.InsertAfter(".0f5")

# 2) "Eggs? (Chests)" -> "Eggs (Chests)"
#    Drop the "?" after "Eggs", splitting the run into "Eggs" / " (Chests)".
#    A temporary bookmark at the split point forces the run break; it is
#    then removed so no bookmark remains on this paragraph (the
#    pre-existing "_GoBack" that used to sit here moves to the Seeds
#    paragraph in step 3 below, matching real Word's "last edit" marker).
$r2 = $d.Content
$r2.Find.Execute("Eggs", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint2 = $d.Range($r2.End, $r2.End)
$d.Bookmarks.Add("TempMarker", $splitPoint2)
$qMark2 = $d.Range($r2.End, $r2.End + 1)
$qMark2.Delete()
$d.Bookmarks("TempMarker").Delete()

# 3) "Seeds? (Coins)" -> "Seeds (Coins)"
#    Same kind of edit as step 2, but this is the last edit made, so the
#    "_GoBack" bookmark ends up sitting right after "Seeds".
$r3 = $d.Content
$r3.Find.Execute("Seeds", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint3 = $d.Range($r3.End, $r3.End)
$qMark3 = $d.Range($r3.End, $r3.End + 1)
$qMark3.Delete()
$d.Bookmarks.Add("_GoBack", $splitPoint3)

# 4) "Goal" -> "Nest (Goal)"
$r4 = $d.Content
$r4.Find.Execute("Goal", $false, $false, $false, $false, $false, $true, 1, $false, "Nest (Goal)", 2)
